$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (46082 -> 46083) for every data row (rows 2 through 95).
for ($r = 2; $r -le 95; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46082) {
        $cell.Value2 = 46083
    }
}
